# Fix #12 - Remove "1-" from Excel tabs
$wb = $excel.ActiveWorkbook

# --- Rename the worksheet tabs, dropping the "1-" prefix ---
$wb.Worksheets.Item("1-Data").Name           = "Data"
$wb.Worksheets.Item("1-Error Frequency").Name = "Error Frequency"
$wb.Worksheets.Item("1-Error Count").Name     = "Error Count"
$wb.Worksheets.Item("1-Histogram").Name       = "Histogram"
$wb.Worksheets.Item("1-Rules").Name           = "Rules"

# --- Repoint the chart series formulas at the renamed sheets ---
# (defined names already auto-update when the sheet is renamed, but the
#  chart series cached formulas reference sheet names directly and need
#  to be re-pointed by hand)

$wsData = $wb.Worksheets.Item("Data")
$chartData1 = $wsData.ChartObjects().Item(1).Chart
$chartData1.SeriesCollection().Item(1).Formula = "=SERIES(Data!`$A`$133,,Data!`$A`$133:`$A`$137,1)"
$chartData1.SeriesCollection().Item(2).Formula = "=SERIES(Data!`$A`$136,,[0]!no_errors_warnings,2)"

$wsErrFreq = $wb.Worksheets.Item("Error Frequency")
$chartErrFreq = $wsErrFreq.ChartObjects().Item(1).Chart
$chartErrFreq.SeriesCollection().Item(1).Formula = "=SERIES(""Errors"",'Error Frequency'!`$A`$2:`$A`$15,'Error Frequency'!`$B`$2:`$B`$15,1)"

$wsErrCount = $wb.Worksheets.Item("Error Count")
$chartErrCount1 = $wsErrCount.ChartObjects().Item(1).Chart
$chartErrCount1.SeriesCollection().Item(1).Formula = "=SERIES(,'Error Count'!`$A`$2:`$A`$56,'Error Count'!`$B`$2:`$B`$56,1)"
$chartErrCount2 = $wsErrCount.ChartObjects().Item(2).Chart
$chartErrCount2.SeriesCollection().Item(1).Formula = "=SERIES(,'Error Count'!`$D`$2:`$D`$52,'Error Count'!`$E`$2:`$E`$52,1)"

$wsHist = $wb.Worksheets.Item("Histogram")
$chartHist1 = $wsHist.ChartObjects().Item(1).Chart
$chartHist1.SeriesCollection().Item(1).Formula = "=SERIES(""Frequency"",Histogram!`$A`$2:`$A`$9,Histogram!`$B`$2:`$B`$9,1)"
$chartHist2 = $wsHist.ChartObjects().Item(2).Chart
$chartHist2.SeriesCollection().Item(1).Formula = "=SERIES(""Frequency"",Histogram!`$A`$14:`$A`$21,Histogram!`$B`$14:`$B`$21,1)"

# --- Fix the active tab / selection: Data tab is now the active / selected one ---
$wsData.Activate()
$wsData.Select()
$wsData.Range("A110").Select()
$excel.ActiveWindow.ScrollRow = 110
$excel.ActiveWindow.ScrollColumn = 1
